# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (interested-attendee count) figures in column F
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 58
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 9588
$ws1.Range("F8").Value = 331
$ws1.Range("F9").Value = 1212
$ws1.Range("F10").Value = 2221
$ws1.Range("F13").Value = 22
$ws1.Range("F15").Value = 455
$ws1.Range("F16").Value = 97
$ws1.Range("F18").Value = 1336

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 58
$ws4.Range("F6").Value = 186
$ws4.Range("F7").Value = 9588
$ws4.Range("F9").Value = 331
$ws4.Range("F10").Value = 1212
$ws4.Range("F11").Value = 2221
$ws4.Range("F14").Value = 22
$ws4.Range("F16").Value = 455
$ws4.Range("F17").Value = 97
$ws4.Range("F19").Value = 1336
